# quarterly_seprated.xlsx update
#
# The reporting window rolled forward by one quarter: the oldest quarter
# column (فصل سوم منتهی به 1399/06, column E) is dropped and a new quarter
# (فصل اول منتهی به 1401/12) is appended as the new last column (N).
# Every other quarter's figures shift one column to the left (F->E, G->F, ...,
# N->M) and the freshly reported quarter's figures land in the new N column.
#
# This applies uniformly to the quarter-header rows (8, 27, 46, 58, 77) and to
# every data row in the five tables (quantities sold, sales amount, unit
# price, cost, gross profit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 5 + $i).Value = $vals[$i]
    }
}

# --- Quarter header rows (columns E:N) ---------------------------------
$quarterLabels = @(
    "فصل چهارم منتهی به 1399/09",
    "فصل اول منتهی به 1399/12",
    "فصل دوم منتهی به 1400/03",
    "فصل سوم منتهی به 1400/06",
    "فصل چهارم منتهی به 1400/09",
    "فصل اول منتهی به 1400/12",
    "فصل دوم منتهی به 1401/03",
    "فصل سوم منتهی به 1401/06",
    "فصل چهارم منتهی به 1401/09",
    "فصل اول منتهی به 1401/12"
)
$headerRows = @(8, 27, 46, 58, 77)
foreach ($hr in $headerRows) {
    Set-RowValues $hr $quarterLabels
}

# --- مقدار فروش داخلی (domestic quantities sold) ------------------------
Set-RowValues 11 @(-181999, 6512, 7951, "-", 8322, 7, 16417, 7487, 8597, 7678)
Set-RowValues 12 @(43, 7, 2, 4, 4, 0, 0, 0, 0, "-")
Set-RowValues 13 @("-", "-", "-", 21394, "-", "-", "-", "-", "-", "-")
Set-RowValues 14 @(-181956, 6519, 7953, 21398, 8326, 7, 16417, 7487, 8597, 7678)

# --- مقدار فروش خارجی (export quantities sold) --------------------------
Set-RowValues 16 @(0, "-", "-", "-", "-", "-", 12, -12, 0, 0)
Set-RowValues 17 @(0, "-", "-", "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 18 @(0, 0, 0, 0, 0, 0, 12, -12, 0, 0)

# --- مقدار فروش درآمد ارائه خدمات (services quantities sold) ------------
Set-RowValues 20 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# --- برگشت از فروش / سایر تخفیفات / جمع (quantities) --------------------
Set-RowValues 21 @(0, 0, -27, "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 22 @(0, 0, 0, 0, 0, "-", 0, 0, 0, 0)
Set-RowValues 23 @(-181956, 6519, 7925, 21398, 8326, 7, 16429, 7475, 8597, 7678)

# --- فروش داخلی (domestic sales amount) ---------------------------------
Set-RowValues 30 @(4280006, 4134589, 5733016, "-", 5463220, 5087108, 6932461, 6513542, 8106599, 7803958)
Set-RowValues 31 @(19960, 3217, 762, 2184, 1809, 0, 0, 1, 0, "-")
Set-RowValues 32 @("-", "-", "-", 14428548, "-", "-", "-", "-", "-", "-")
Set-RowValues 33 @(4299966, 4137806, 5733778, 14430732, 5465029, 5087108, 6932461, 6513543, 8106599, 7803958)

# --- فروش خارجی (export sales amount) -----------------------------------
Set-RowValues 35 @(0, "-", "-", "-", "-", "-", 25057, 54297, -79354, 0)
Set-RowValues 36 @(0, "-", "-", "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 37 @(0, 0, 0, 0, 0, 0, 25057, 54297, -79354, 0)

# --- فروش درآمد ارائه خدمات (services sales amount) ---------------------
Set-RowValues 39 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# --- برگشت از فروش / سایر تخفیفات / جمع (sales amount) ------------------
Set-RowValues 40 @(0, 0, -11860, "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 41 @(0, 0, -203, 0, 0, "-", 0, 0, 0, 0)
Set-RowValues 42 @(4299966, 4137806, 5721715, 14430732, 5465029, 5087108, 6957518, 6567840, 8027245, 7803958)

# --- نرخ فروش داخلی (domestic unit price) -------------------------------
Set-RowValues 49 @(-23516646, 634918458, 721026161, "-", 656479212, 686797353854, 422283860, 869980232, 942956729, 1016411672)
Set-RowValues 50 @(464186047, 459571429, 482278481, 546000000, 452250000, "-", "-", "-", "-", "-")
Set-RowValues 51 @("-", "-", "-", 674420305, "-", "-", "-", "-", "-", "-")

# --- نرخ فروش خارجی (export unit price) ---------------------------------
Set-RowValues 53 @("-", "-", "-", "-", "-", "-", 2088083333, -4524750000, "-", "-")
Set-RowValues 54 @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")

# --- بهای تمام شده داخلی (domestic cost) ---------------------------------
Set-RowValues 61 @(-2585669, -2706033, -3973297, "-", -14621235, -3989842, -6242002, -5414506, -6246488, -5698607)
Set-RowValues 62 @(-9473, -5293, 2983, -1190, -996, 0, 0, -1, 0, "-")
Set-RowValues 63 @("-", "-", "-", -10207518, "-", "-", "-", "-", "-", "-")
Set-RowValues 64 @(-2595142, -2711326, -3970314, -10208708, -14622231, -3989842, -6242002, -5414507, -6246488, -5698607)

# --- بهای تمام شده خارجی (export cost) -----------------------------------
Set-RowValues 66 @(0, "-", "-", "-", "-", "-", -19864, -72607, 92471, 0)
Set-RowValues 67 @(0, "-", "-", "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 68 @(0, 0, 0, 0, 0, 0, -19864, -72607, 92471, 0)

# --- بهای تمام شده درآمد ارائه خدمات (services cost) ---------------------
Set-RowValues 70 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# --- برگشت از فروش / سایر تخفیفات / جمع (cost) ---------------------------
Set-RowValues 71 @(0, 0, 11842, "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 72 @(0, 0, 0, 0, 0, "-", 0, 0, 0, 0)
Set-RowValues 73 @(-2595142, -2711326, -3958472, -10208708, -14622231, -3989842, -6261866, -5487114, -6154017, -5698607)

# --- سود ناخالص داخلی (domestic gross profit) -----------------------------
Set-RowValues 80 @(1694337, 1428556, 1759719, "-", 5270533, 1097266, 690459, 1099036, 1860111, 2105351)
Set-RowValues 81 @(10487, -2076, 3745, 994, 813, 0, 0, 0, 0, "-")
Set-RowValues 82 @("-", "-", "-", 4221030, "-", "-", "-", "-", "-", "-")
Set-RowValues 83 @(1704824, 1426480, 1763464, 4222024, 5271346, 1097266, 690459, 1099036, 1860111, 2105351)

# --- سود ناخالص خارجی (export gross profit) -------------------------------
Set-RowValues 85 @(0, "-", "-", "-", "-", "-", 5193, -18310, 13117, 0)
Set-RowValues 86 @(0, "-", "-", "-", "-", "-", "-", "-", "-", "-")
Set-RowValues 87 @(0, 0, 0, 0, 0, 0, 5193, -18310, 13117, 0)

# --- سود ناخالص درآمد ارائه خدمات (services gross profit) -----------------
Set-RowValues 89 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# --- جمع (gross profit total) ---------------------------------------------
Set-RowValues 90 @(1704824, 1426480, 1763464, 4222024, 5271346, 1097266, 695652, 1080726, 1873228, 2105351)
